$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.211.09'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '3.326.64'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.91'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.322.03'
$ws.Range("E8").Value = '  +0.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.578'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.79%  '

$ws.Range("E10").Value = '  -2.16%  '

$ws.Range("E11").Value = '  -1.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.10'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("E13").Value = '  -1.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '668.08'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +10.23%  '

$ws.Range("D15").Value = '3.860.41'
$ws.Range("E15").Value = '  +0.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.48'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.35%  '

$ws.Range("D17").Value = '66.254.04'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.92'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.18%  '

$ws.Range("E19").Value = '  -0.30%  '

$ws.Range("D20").Value = '3.325.73'
$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.15'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.896'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.10%  '

$ws.Range("E23").Value = '  -3.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.63'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("E25").Value = '  -1.72%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.53'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '31.80'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.49'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.80'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '600.27'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +7.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.90'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.92%  '

$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.846.65'
$ws.Range("E35").Value = '  +3.43%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.105'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.05'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("E40").Value = '  -4.07%  '

$ws.Range("E41").Value = '  -2.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.42'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.16'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.59%  '

$ws.Range("E46").Value = '  -2.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -12.14%  '

$ws.Range("E48").Value = '  -1.77%  '

$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.39%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.55'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.31'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.11%  '
